# The data table previously had 4 data rows (Sending x Target cluster
# combinations of ECs/FAPs -> Hc/C5ar1). The underlying TPM data was
# recomputed and the workbook now only reports a single row: the
# FAPs -> Hc/C5ar1 edge (with specificity columns renormalised to 1
# since it is now the only row). Remove the other three data rows and
# refresh the specificity-derived columns on the remaining row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 2-4 (old ECs->ECs, ECs->FAPs, FAPs->ECs rows); the former
# row 5 (FAPs->FAPs / FAPs->Hc/C5ar1->FAPs) shifts up to become row 2.
$ws.Rows("2:4").Delete()

# Specificity values are recomputed now that FAPs->Hc/C5ar1 is the
# sole remaining edge, so each specificity column normalises to 1.
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
